$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: style a single-run paragraph (Heading1 + Arial/Bold/19pt) and set
# its text, matching the OOXML shape produced in the target diff:
#   <w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr>
#     <w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:b/>
#           <w:sz w:val="38"/></w:rPr><w:t>...</w:t></w:r></w:p>
# ---------------------------------------------------------------------------
function Style-FieldParagraph($para, [string]$text) {
    $para.Style = "Heading1"
    $rng = $para.Range
    $rng.End = $rng.End - 1
    $rng.Text = $text
    $rng.Font.Name = "Arial"
    $rng.Font.Bold = $true
    $rng.Font.Size = 19
}

# Turn a (currently blank) paragraph into the thin "spacing" spacer paragraph:
#   <w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto" w:after="0"/></w:pPr></w:p>
function Style-SpacerParagraph($para) {
    $para.LineSpacingRule = 0
    $para.SpaceAfter = 0
}

# Turn a (currently blank) paragraph into a page-break paragraph:
#   <w:p><w:r><w:br w:type="page"/></w:r></w:p>
function Style-PageBreakParagraph($para) {
    $rng = $para.Range
    $rng.End = $rng.End - 1
    $rng.Text = [char]12
}

# ---------------------------------------------------------------------------
# Record 1 (was "Муратов Григорий Степанович" / paragraphs 1-5): prepend the
# (still-blank) spacer paragraph *before* any styling happens, so it does
# not inherit the Heading1/Arial formatting applied afterwards, then retext
# + reformat the 5 field paragraphs (now shifted to 2-6).
# ---------------------------------------------------------------------------
$d.Paragraphs(1).Range.InsertParagraphBefore()
Style-SpacerParagraph $d.Paragraphs(1)

Style-FieldParagraph $d.Paragraphs(2) "Хооп Антон Петрович"
Style-FieldParagraph $d.Paragraphs(3) "26 лет"
Style-FieldParagraph $d.Paragraphs(4) "Родился в п. Сопотно Лядского р-на Ленинградской обл."
Style-FieldParagraph $d.Paragraphs(5) "зав. Вагошской мельницей"
Style-FieldParagraph $d.Paragraphs(6) "Расстрелян 08 июля 1938 в г. Ленинград"

# Paragraph 7 is now the original page break after record 1 - leave as is.

# ---------------------------------------------------------------------------
# Record 2 (was "Мусатов Сергей Никанорович" / originally paragraphs 7-11,
# now shifted to 8-12 after the spacer insertion above): same approach -
# spacer first, then retext + reformat the 5 field paragraphs.
# ---------------------------------------------------------------------------
$d.Paragraphs(8).Range.InsertParagraphBefore()
Style-SpacerParagraph $d.Paragraphs(8)

Style-FieldParagraph $d.Paragraphs(9) "Худяков Кузьма Авдеевич"
Style-FieldParagraph $d.Paragraphs(10) "48 лет"
Style-FieldParagraph $d.Paragraphs(11) "Родился в Томская обл., Чаинский р-н, дер. Красноярка"
Style-FieldParagraph $d.Paragraphs(12) "единоличник"
Style-FieldParagraph $d.Paragraphs(13) "Расстрелян 31 марта 1938 в неизвестно"

# Paragraph 14 is now the original page break after record 2 - leave as is.

# ---------------------------------------------------------------------------
# Three brand-new records appended at the end of the document, each as:
#   spacer paragraph, 5 styled field paragraphs, page-break paragraph.
# ---------------------------------------------------------------------------
$newRecords = @(
    @{
        name = "Чекулаев Михаил Васильевич"
        age  = "41 лет"
        born = "Родился в с. Брусяны"
        job  = "колхозник"
        ex   = "Расстрелян 11 мая 1938 в Куйбышев"
    },
    @{
        name = "Шестериков Николай Васильевич"
        age  = "25 лет"
        born = "Родился в в с. Ставрополье."
        job  = "радист."
        ex   = "Расстрелян 12 декабря 1941 в в Куйбышеве."
    },
    @{
        name = "Шишкова Михалина Иосифовна"
        age  = "54 лет"
        born = "Родилась в дер. Залесье Ошмянского уезд Виленской губ."
        job  = "Домохозяйка"
        ex   = "Расстреляна 15 января 1938 в г. Ленинград"
    }
)

foreach ($rec in $newRecords) {
    $anchor = $d.Paragraphs($d.Paragraphs.Count).Range
    $anchor.InsertParagraphAfter()
    $anchor.InsertParagraphAfter()
    $anchor.InsertParagraphAfter()
    $anchor.InsertParagraphAfter()
    $anchor.InsertParagraphAfter()
    $anchor.InsertParagraphAfter()
    $anchor.InsertParagraphAfter()

    $base = $d.Paragraphs.Count - 7

    Style-SpacerParagraph $d.Paragraphs($base + 1)
    Style-FieldParagraph $d.Paragraphs($base + 2) $rec.name
    Style-FieldParagraph $d.Paragraphs($base + 3) $rec.age
    Style-FieldParagraph $d.Paragraphs($base + 4) $rec.born
    Style-FieldParagraph $d.Paragraphs($base + 5) $rec.job
    Style-FieldParagraph $d.Paragraphs($base + 6) $rec.ex
    Style-PageBreakParagraph $d.Paragraphs($base + 7)
}

Write-Output "done"
